$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Narrow column A from 20 to 17 (stored XML width = ColumnWidth + 0.8333...)
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668

# Update existing Threshold row (row 2) BTC/ETH values
$ws.Cells.Item(2, 3).Value = 0.0003182827846920551
$ws.Cells.Item(2, 4).Value = 0.004508441233855234

# Update existing timestamp row (row 3) with the new snapshot
$ws.Cells.Item(3, 1).Value = "24 Mar 22, 17:37PM"
$ws.Cells.Item(3, 2).Value = 13.98
$ws.Cells.Item(3, 3).Value = 0.0003182880746132817
$ws.Cells.Item(3, 4).Value = 0.004509815623382007

# Append a new row 4 with the same snapshot
$ws.Cells.Item(4, 1).Value = "24 Mar 22, 17:37PM"
$ws.Cells.Item(4, 2).Value = 13.98
$ws.Cells.Item(4, 3).Value = 0.0003182880746132817
$ws.Cells.Item(4, 4).Value = 0.004509815623382007

# Append a new row 5 with the latest snapshot
$ws.Cells.Item(5, 1).Value = "24 Mar 22, 17:37PM"
$ws.Cells.Item(5, 2).Value = 13.98
$ws.Cells.Item(5, 3).Value = 0.0003181687680484818
$ws.Cells.Item(5, 4).Value = 0.004509815623382007
